$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "1" to "C.Rustavi"
$ws.Name = "C.Rustavi"

# The "Rural" row (row 7) values become confidential/unavailable, shown as "..."
$ws.Range("B7:O7").Value = "..."

# Add a footnote in row 8 explaining the "..." marker, with "Note:" in bold+underline
# Row 8 reuses the same (borderless, fill-less, small-font) style already present on row 9
$ws.Range("A9").Copy()
$ws.Range("A8").PasteSpecial(-4122)  # xlPasteFormats

$noteLabel = "Note:"
$noteBody = " „ ... „ - Data is confidential or unavailable."
$noteText = $noteLabel + $noteBody

$noteCell = $ws.Range("A8")
$noteCell.Value = $noteText

$labelChars = $noteCell.Characters(1, $noteLabel.Length)
$labelFont = $labelChars.Font
$labelFont.Bold = $true
$labelFont.Underline = $true
$labelFont.Size = 9
$labelFont.Name = "Arial"

$bodyChars = $noteCell.Characters($noteLabel.Length + 1, $noteBody.Length)
$bodyFont = $bodyChars.Font
$bodyFont.Size = 9
$bodyFont.Name = "Arial"
